$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New loss-of-sale records (columns A-K) to append below the existing data (row 36)
$data = @(
    @("35","22-12-2025","ivin","9020244484","03-01-2026","Noel Jacob","Loss","CUSTOMER INTERNAL ISSUES","FAMILY DISAPPROVEL","-","They are looking for a suit,we are showed multiple items ,he will confirm with family and will revisit again"),
    @("36","22-12-2025","saji","9544577340","28-12-2025","Noel Jacob","Loss","ENQUIRY","Enquiry for Relative/Friend","-","they need green colour non premium 3peace suit  different sizes. and they want opinion from cousins revist again"),
    @("37","22-12-2025","febin","9207420831","12-01-2026","ABHIJITH KUMAR P A","Loss","ENQUIRY","ENQUIRY WITHOUT BRIDE/FAMILY","-","they have checked the products which is suits for him,will revisit with bride in upcoming days"),
    @("38","22-12-2025","jims","9731146756","15-01-2026","Noel Jacob","Loss","ENQUIRY","ENQUIRY WITHOUT BRIDE/FAMILY","-","just enquiry revist on 3 days"),
    @("39","22-12-2025","Hari","9747511995","29-12-2025","Jithinsha R","Loss","ENQUIRY","ENQUIRY WITHOUT BRIDE/FAMILY","-","Customer ok with Navy blue and Peacock blue bengala,He needs a suires for wedding also.he will revisit the store by tomorrw for booking for both items"),
    @("40","23-12-2025","Alvin","9495564724","27-12-2025","Aswin Raj M. R","Loss","ENQUIRY","ENQUIRY WITHOUT BRIDE/FAMILY","-","want black embro but dates are not available"),
    @("41","23-12-2025","sachin","9495574128","27-12-2025","Arjun Reji","Loss","CUSTOMER INTERNAL ISSUES","FAMILY DISAPPROVEL","-","He was okey with the product but family was disapprovel"),
    @("42","23-12-2025","Ashok","8714298728","07-01-2026","Aswin Raj M. R","Loss","ENQUIRY","ENQUIRY WITHOUT BRIDE/FAMILY","-","Customer is ok with biege suit he discuss with bride and family and visit again"),
    @("43","24-12-2025","rasal","9072226848","11-01-2026","Arjun Reji","Loss","ENQUIRY","ENQUIRY WITHOUT TRIAL","-","Enquiry and price details"),
    @("44","24-12-2025","Deepak","8848355981","28-12-2025","Aswin Raj M. R","Loss","ENQUIRY","ENQUIRY WITHOUT TRIAL","-","just visit"),
    @("45","24-12-2025","johnson","8891512403","06-01-2026","Noel Jacob","Loss","ENQUIRY","ENQUIRY WITHOUT TRIAL","-","want black colour short work model  and that product already booked for another custmor"),
    @("46","24-12-2025","philip","7736254684","08-01-2026","Arjun Reji","Loss","PRODUCT","PRODUCT NOT AVAILABLE","-","He was ok with the product.but the product not for available at the date"),
)

$startRow = 37
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $startRow + $idx
    $prevRow = $row - 1
    $rec = $data[$idx]

    # Column A (#) - numeric; copy number format from the row above to keep the "0" integer format
    $ws.Cells.Item($row, 1).Value = [int]$rec[0]
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    # Column B (Date) - plain text; leading apostrophe keeps it literal instead of an auto-parsed date
    $ws.Cells.Item($row, 2).Value = "'" + $rec[1]
    $ws.Cells.Item($row, 2).Style = $ws.Cells.Item($prevRow, 2).Style

    # Column C (Customer Name)
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item($prevRow, 3).Style

    # Column D (Contact) - numeric; copy number format from the row above to keep the "0" integer format
    $ws.Cells.Item($row, 4).Value = [double]$rec[3]
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($prevRow, 4).NumberFormat

    # Column E (Function Date) - plain text; leading apostrophe keeps it literal instead of an auto-parsed date
    $ws.Cells.Item($row, 5).Value = "'" + $rec[4]
    $ws.Cells.Item($row, 5).Style = $ws.Cells.Item($prevRow, 5).Style

    # Column F (Staff)
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 6).Style = $ws.Cells.Item($prevRow, 6).Style

    # Column G (Status)
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 7).Style = $ws.Cells.Item($prevRow, 7).Style

    # Column H (Category)
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 8).Style = $ws.Cells.Item($prevRow, 8).Style

    # Column I (Sub Category)
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 9).Style = $ws.Cells.Item($prevRow, 9).Style

    # Column J (Repeat count)
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 10).Style = $ws.Cells.Item($prevRow, 10).Style

    # Column K (Remarks)
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $ws.Cells.Item($row, 11).Style = $ws.Cells.Item($prevRow, 11).Style
}

$wb.Save()
